$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '57.144.39'
$ws.Cells.Item(2, 5).Value = '  +0.63%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.427.11'
$ws.Cells.Item(3, 5).Value = '  -1.67%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '489.20'
$ws.Cells.Item(5, 5).Value = '  -0.15%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '154.44'
$ws.Cells.Item(6, 5).Value = '  +1.79%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.616'
$ws.Cells.Item(7, 5).Value = '  +19.67%  '
$ws.Cells.Item(8, 5).Value = '  -0.12%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2.447.33'
$ws.Cells.Item(9, 5).Value = '  -1.16%  '
$ws.Cells.Item(10, 5).Value = '  +8.64%  '
$ws.Cells.Item(11, 5).Value = '  +0.70%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.333'
$ws.Cells.Item(12, 5).Value = '  +0.02%  '
$ws.Cells.Item(13, 5).Value = '  +1.16%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '2.855.44'
$ws.Cells.Item(14, 5).Value = '  -1.59%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '57.144.58'
$ws.Cells.Item(15, 5).Value = '  +0.10%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '20.60'
$ws.Cells.Item(16, 5).Value = '  -1.82%  '
$ws.Cells.Item(17, 5).Value = '  -2.67%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.442.06'
$ws.Cells.Item(18, 5).Value = '  -0.93%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.62'
$ws.Cells.Item(19, 5).Value = '  +1.70%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '324.32'
$ws.Cells.Item(20, 5).Value = '  +1.26%  '
$ws.Cells.Item(21, 5).Value = '  -1.22%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.00'
$ws.Cells.Item(22, 5).Value = '  +0.29%  '
$ws.Cells.Item(23, 5).Value = '  +1.38%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '57.85'
$ws.Cells.Item(24, 5).Value = '  -0.37%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.995'
$ws.Cells.Item(25, 5).Value = '  -0.64%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.400'
$ws.Cells.Item(26, 5).Value = '  -1.21%  '
$ws.Cells.Item(27, 5).Value = '  -1.62%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.532.93'
$ws.Cells.Item(28, 5).Value = '  -2.05%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.26'
$ws.Cells.Item(29, 5).Value = '  -4.13%  '
$ws.Cells.Item(30, 5).Value = '  -2.15%  '
$ws.Cells.Item(31, 5).Value = '  -0.02%  '
$ws.Cells.Item(32, 2).Value = 'EthereumClassic'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '18.70'
$ws.Cells.Item(32, 5).Value = '  +2.46%  '
$ws.Cells.Item(33, 2).Value = 'Monero'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '150.22'
$ws.Cells.Item(33, 5).Value = '  -0.33%  '
$ws.Cells.Item(34, 5).Value = '  +0.33%  '
$ws.Cells.Item(35, 5).Value = '  +1.77%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '3.78'
$ws.Cells.Item(36, 5).Value = '  +0.79%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.13'
$ws.Cells.Item(37, 5).Value = '  -1.38%  '
$ws.Cells.Item(38, 5).Value = '  -7.72%  '
$ws.Cells.Item(39, 2).Value = 'Bittensor'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '285.38'
$ws.Cells.Item(39, 5).Value = '  +8.82%  '
$ws.Cells.Item(40, 2).Value = 'Stellar'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.101'
$ws.Cells.Item(40, 5).Value = '  +6.98%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '34.05'
$ws.Cells.Item(41, 5).Value = '  -0.20%  '
$ws.Cells.Item(42, 2).Value = 'Stacks'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.38'
$ws.Cells.Item(42, 5).Value = '  -0.68%  '
$ws.Cells.Item(43, 2).Value = 'Filecoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.52'
$ws.Cells.Item(43, 5).Value = '  +0.63%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.994'
$ws.Cells.Item(44, 5).Value = '  -0.21%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.603'
$ws.Cells.Item(45, 5).Value = '  -0.66%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0532'
$ws.Cells.Item(46, 5).Value = '  -4.28%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '10.21'
$ws.Cells.Item(47, 5).Value = '  +0.02%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0227'
$ws.Cells.Item(48, 5).Value = '  -0.35%  '
$ws.Cells.Item(49, 5).Value = '  -4.20%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.902.57'
$ws.Cells.Item(50, 5).Value = '  +2.89%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '17.60'
$ws.Cells.Item(51, 5).Value = '  -0.40%  '
